$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 112.875
$ws.Range("I11").Value = 112.875
$ws.Range("K11").Value = 112.875
$ws.Range("M11").Value = 27.125

$ws.Range("H131").Value = 2166.182
$ws.Range("I131").Value = 648.1111
$ws.Range("K131").Value = 1944.3333
$ws.Range("M131").Value = 3095.6667

$ws.Range("H141").Value = 1656.45
$ws.Range("I141").Value = 1300.2222
$ws.Range("K141").Value = 3900.6666
$ws.Range("M141").Value = 1279.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5996.6816
$ws.Range("I32").Value = 4485.488
$ws.Range("J32").Value = 26649.666
$ws.Range("K32").Value = 4485.488
$ws.Range("L32").Value = 26649.666
$ws.Range("M32").Value = -4198.488
$ws.Range("N32").Value = -27223.666

$ws.Range("H61").Value = 13425.586
$ws.Range("I61").Value = 2874.5
$ws.Range("J61").Value = 36872.445
$ws.Range("K61").Value = 2874.5
$ws.Range("L61").Value = 36872.445
$ws.Range("M61").Value = -2662.5
$ws.Range("N61").Value = -37296.445

$ws.Range("H132").Value = 3721291.8
$ws.Range("I132").Value = 7165.25
$ws.Range("J132").Value = 6692593
$ws.Range("K132").Value = 21495.75
$ws.Range("L132").Value = 20077779
$ws.Range("M132").Value = -18965.75
$ws.Range("N132").Value = -20082839

$ws.Range("H136").Value = 13425.586
$ws.Range("I136").Value = 2874.5
$ws.Range("J136").Value = 36872.445
$ws.Range("K136").Value = 8623.5
$ws.Range("L136").Value = 110617.335
$ws.Range("M136").Value = -6073.5
$ws.Range("N136").Value = -115717.335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 13653.77
$ws.Range("I20").Value = 4391.5
$ws.Range("J20").Value = 25151.758
$ws.Range("K20").Value = 4391.5
$ws.Range("L20").Value = 25151.758
$ws.Range("M20").Value = -4144.5
$ws.Range("N20").Value = -25645.758

$ws.Range("H99").Value = 12866.333
$ws.Range("I99").Value = 1371.8889
$ws.Range("J99").Value = 24360.777
$ws.Range("K99").Value = 1371.8889
$ws.Range("L99").Value = 24360.777
$ws.Range("M99").Value = 126.1111000000001
$ws.Range("N99").Value = -27356.777

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

$ws.Range("H134").Value = 7962.8613
$ws.Range("I134").Value = 2026.4333
$ws.Range("J134").Value = 37645
$ws.Range("K134").Value = 6079.2999
$ws.Range("L134").Value = 112935
$ws.Range("M134").Value = -3544.2999
$ws.Range("N134").Value = -118005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()

$ws.Range("H31").Value = 18185.156
$ws.Range("I31").Value = 8162.7646
$ws.Range("K31").Value = 8162.7646
$ws.Range("M31").Value = -7867.7646

$ws.Range("H34").Value = 18185.156
$ws.Range("I34").Value = 8162.7646
$ws.Range("K34").Value = 8162.7646
$ws.Range("M34").Value = -7960.7646

$ws.Range("H58").Value = 13709.771
$ws.Range("I58").Value = 8163.6665
$ws.Range("J58").Value = 15629.577
$ws.Range("K58").Value = 8163.6665
$ws.Range("L58").Value = 15629.577
$ws.Range("M58").Value = -7960.6665
$ws.Range("N58").Value = -16035.577

$ws.Range("H105").Value = 22319
$ws.Range("I105").Value = 25558.5
$ws.Range("K105").Value = 25558.5
$ws.Range("M105").Value = -23811.5

$ws.Range("H110").Value = 79398.8
$ws.Range("J110").Value = 79398.8
$ws.Range("L110").Value = 79398.8
$ws.Range("N110").Value = -87578.8

$ws.Range("H122").Value = 4175.032
$ws.Range("I122").Value = 2099.7896
$ws.Range("J122").Value = 7460.8335
$ws.Range("K122").Value = 6299.3688
$ws.Range("L122").Value = 22382.5005
$ws.Range("M122").Value = -3849.3688
$ws.Range("N122").Value = -27282.5005

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws.Range("H132").Value = 6149.7407
$ws.Range("I132").Value = 2597.5833
$ws.Range("J132").Value = 8991.467000000001
$ws.Range("K132").Value = 7792.749899999999
$ws.Range("L132").Value = 26974.401
$ws.Range("M132").Value = -5262.749899999999
$ws.Range("N132").Value = -32034.401

$ws.Range("H134").Value = 33341100
$ws.Range("I134").Value = 1856.4615
$ws.Range("K134").Value = 5569.3845
$ws.Range("M134").Value = -3034.3845

$ws.Range("H136").Value = 13709.771
$ws.Range("I136").Value = 8163.6665
$ws.Range("J136").Value = 15629.577
$ws.Range("K136").Value = 24490.9995
$ws.Range("L136").Value = 46888.731
$ws.Range("M136").Value = -21940.9995
$ws.Range("N136").Value = -51988.731

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 12999.667
$ws.Range("I69").Value = 3999
$ws.Range("K69").Value = 11997
$ws.Range("M69").Value = -11186

$ws.Range("H72").Value = 12999.667
$ws.Range("I72").Value = 3999
$ws.Range("K72").Value = 35991
$ws.Range("M72").Value = -31935

$ws.Range("H74").Value = 23125
$ws.Range("J74").Value = 23125
$ws.Range("L74").Value = 69375
$ws.Range("N74").Value = -71497

$ws.Range("H76").Value = 8417.857
$ws.Range("I76").Value = 4820.8335
$ws.Range("K76").Value = 14462.5005
$ws.Range("M76").Value = -14079.5005

$ws.Range("H77").Value = 23125
$ws.Range("J77").Value = 23125
$ws.Range("L77").Value = 208125
$ws.Range("N77").Value = -218733

$ws.Range("H79").Value = 8417.857
$ws.Range("I79").Value = 4820.8335
$ws.Range("K79").Value = 14462.5005
$ws.Range("M79").Value = -13136.5005

$ws.Range("H80").Value = 14321
$ws.Range("I80").Value = 11199
$ws.Range("J80").Value = 16055.444
$ws.Range("K80").Value = 33597
$ws.Range("L80").Value = 48166.33199999999
$ws.Range("M80").Value = -32661
$ws.Range("N80").Value = -50038.33199999999

$ws.Range("H83").Value = 14321
$ws.Range("I83").Value = 11199
$ws.Range("J83").Value = 16055.444
$ws.Range("K83").Value = 100791
$ws.Range("L83").Value = 144498.996
$ws.Range("M83").Value = -96111
$ws.Range("N83").Value = -153858.996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 10263.333
$ws.Range("I20").Value = 7895
$ws.Range("K20").Value = 7895
$ws.Range("M20").Value = -7650

$ws.Range("H24").Value = 10000
$ws.Range("I24").Value = 5000
$ws.Range("K24").Value = 5000
$ws.Range("M24").Value = -4827

$ws.Range("H26").Value = 24862.955
$ws.Range("I26").Value = 10000
$ws.Range("K26").Value = 10000
$ws.Range("M26").Value = -9720

$ws.Range("H50").Value = 24862.955
$ws.Range("I50").Value = 10000
$ws.Range("K50").Value = 10000
$ws.Range("M50").Value = -9502

$ws.Range("H70").Value = 18493.375
$ws.Range("J70").Value = 24689
$ws.Range("L70").Value = 24689
$ws.Range("N70").Value = -25229

$ws.Range("H73").Value = 18493.375
$ws.Range("J73").Value = 24689
$ws.Range("L73").Value = 24689
$ws.Range("N73").Value = -26561

$ws.Range("H80").Value = 14998.117
$ws.Range("I80").Value = 8180.1113
$ws.Range("K80").Value = 8180.1113
$ws.Range("M80").Value = -7182.1113

$ws.Range("H83").Value = 14998.117
$ws.Range("I83").Value = 8180.1113
$ws.Range("K83").Value = 40900.5565
$ws.Range("M83").Value = -35908.5565

$ws.Range("H107").Value = 475.6875
$ws.Range("J107").Value = 526.5454999999999
$ws.Range("L107").Value = 526.5454999999999
$ws.Range("N107").Value = -4366.5455

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H122").Value = 3688.0908
$ws.Range("I122").Value = 3574.3333
$ws.Range("J122").Value = 4200
$ws.Range("K122").Value = 10722.9999
$ws.Range("L122").Value = 12600
$ws.Range("M122").Value = -8272.999899999999
$ws.Range("N122").Value = -17500

$ws.Range("H132").Value = 23929.75
$ws.Range("I132").Value = 19327.75
$ws.Range("K132").Value = 57983.25
$ws.Range("M132").Value = -55453.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 2334.9143
$ws.Range("I55").Value = 1182.8823
$ws.Range("J55").Value = 3422.9443
$ws.Range("K55").Value = 1182.8823
$ws.Range("L55").Value = 3422.9443
$ws.Range("M55").Value = -1009.8823
$ws.Range("N55").Value = -3768.9443

$ws.Range("H132").Value = 3353830.5
$ws.Range("I132").Value = 4140.143
$ws.Range("J132").Value = 8043397
$ws.Range("K132").Value = 12420.429
$ws.Range("L132").Value = 24130191
$ws.Range("M132").Value = -9890.429
$ws.Range("N132").Value = -24135251

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 5000
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()

$ws.Range("H81").Value = 2409.3333
$ws.Range("I81").Value = 1778.4
$ws.Range("K81").Value = 3556.8
$ws.Range("M81").Value = -2495.8

$ws.Range("H84").Value = 2409.3333
$ws.Range("I84").Value = 1778.4
$ws.Range("K84").Value = 17784
$ws.Range("M84").Value = -12480

$ws.Range("H132").Value = 11777.903
$ws.Range("I132").Value = 7605.3184
$ws.Range("K132").Value = 22815.9552
$ws.Range("M132").Value = -20285.9552
